{"js": "// The document contains four `<div>` blocks, each with an `<id>` element\n// such as `<id>p055r_1</id>`. In the original file each of these was split\n// across three separate runs (one for the literal \"<id>\", one for the\n// numeric id itself, one for the literal \"</id>\") using two different\n// run-formatting styles. The edit collapses each of those three runs into\n// a single run containing the full \"<id>...</id>\" text (keeping the\n// Courier-New / brownish \"tag\" formatting of the first run), for every\n// one of the four ids (p055r_1 .. p055r_4).\n\nconst ids = [\"p055r_1\", \"p055r_2\", \"p055r_3\", \"p055r_4\"];\n\nfor (const id of ids) {\n  const needle = `<id>${id}</id>`;\n  const searchResults = context.document.body.search(needle, { matchCase: true });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < searchResults.items.length; i++) {\n    // Replacing the whole matched range (which spans the 3 original runs)\n    // with its own text merges it into a single run that inherits the\n    // formatting of the first of those runs - exactly mirroring the diff.\n    searchResults.items[i].insertText(needle, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains four `<div>` blocks, each with an `<id>` element\n# such as `<id>p055r_1</id>`. In the original file each of these was split\n# across three separate runs (one for the literal \"<id>\", one for the\n# numeric id itself, one for the literal \"</id>\") using two different\n# run-formatting styles. The edit collapses each of those three runs into\n# a single run containing the full \"<id>...</id>\" text (keeping the\n# Courier-New / brownish \"tag\" formatting of the first run), for every\n# one of the four ids (p055r_1 .. p055r_4).\n\n$d = $word.ActiveDocument\n\nfor ($n = 1; $n -le 4; $n++) {\n    $needle = \"<id>p055r_$n</id>\"\n\n    $full = $d.Content\n    $found = $full.Find.Execute($needle)\n    if (-not $found) { continue }\n\n    $start = $full.Start\n    $end = $full.End\n\n    # Keep the first 4 characters (\"<id>\") as-is - that run already has the\n    # formatting (Courier New / brown) that the merged run should end up\n    # with. Grab everything after it (\"p055r_N</id>\") as plain text, delete\n    # that range (which removes the two extra runs), then append the text\n    # back onto the first run so it all becomes one run.\n    $tagStart = $d.Range($start, $start + 4)\n    $remainder = $d.Range($start + 4, $end)\n    $remainderText = $remainder.Text\n\n    $remainder.Delete()\n    $tagStart.InsertAfter($remainderText)\n}\n"}
